# Update the lattice-multiplication exercise table: replace the problem
# (top line), the two split-digit line, and the two product-digit lines
# in every cell of the 5x3 table with the new values from the rebuilt
# worksheet, while leaving the "  ----" separator line and the overall
# cell/run/line-break structure untouched.

$d = $word.ActiveDocument
$vt = [char]11   # vertical-tab char used for manual line breaks (w:br) inside Range.Text

# New content for each of the 15 cells, in row-major order (row1col1,
# row1col2, row1col3, row2col1, ...), matching the order Table.Range.Cells
# enumerates them.
$newCellLines = @(
  @("17 x 31","  3    1","  ----","1|    |","7|    |"),
  @("20 x 44","  4    4","  ----","2|    |","0|    |"),
  @("37 x 89","  8    9","  ----","3|    |","7|    |"),
  @("80 x 81","  8    1","  ----","8|    |","0|    |"),
  @("25 x 99","  9    9","  ----","2|    |","5|    |"),
  @("60 x 33","  3    3","  ----","6|    |","0|    |"),
  @("39 x 13","  1    3","  ----","3|    |","9|    |"),
  @("62 x 69","  6    9","  ----","6|    |","2|    |"),
  @("60 x 71","  7    1","  ----","6|    |","0|    |"),
  @("45 x 41","  4    1","  ----","4|    |","5|    |"),
  @("58 x 26","  2    6","  ----","5|    |","8|    |"),
  @("77 x 32","  3    2","  ----","7|    |","7|    |"),
  @("56 x 45","  4    5","  ----","5|    |","6|    |"),
  @("99 x 44","  4    4","  ----","9|    |","9|    |"),
  @("55 x 17","  1    7","  ----","5|    |","5|    |")
)

$table = $d.Tables.Item(1)

$i = 0
foreach ($cell in $table.Range.Cells) {
    $lines = $newCellLines[$i]
    $newText = [string]::Join($vt, $lines)

    $r = $cell.Range
    $null = $r.MoveEnd(1, -1)   # drop the trailing end-of-cell mark from the range
    $r.Text = $newText

    $i = $i + 1
}
